$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 135.3425948491774
$ws.Range("C3").Value = 12.54490927756532
$ws.Range("C4").Value = 12.83867099854853
$ws.Range("C5").Value = 17.38102768675918
$ws.Range("C6").Value = 22.29492876243445
$ws.Range("C7").Value = 7.370474027753285
$ws.Range("C8").Value = 8.352197002761407
$ws.Range("C9").Value = 22.26925293078039
$ws.Range("C10").Value = 38.3717752354713
$ws.Range("C11").Value = 10.68945285479997
$ws.Range("C12").Value = 2.885510375004642
$ws.Range("C13").Value = 6.525437097727064
$ws.Range("C14").Value = 1.910583943669653
$ws.Range("C15").Value = 2.978396471870794
$ws.Range("C16").Value = 19.08016360504247
$ws.Range("C17").Value = 19.65938016029726
$ws.Range("C18").Value = 19.31351160448671
$ws.Range("C19").Value = 6.167485797608717
$ws.Range("C20").Value = 24.45320896441384
$ws.Range("C21").Value = 68.98491828230149
$ws.Range("C22").Value = 10.76043897760825
$ws.Range("C23").Value = 2.294211075442057
$ws.Range("C24").Value = 22.20279783708753
$ws.Range("C25").Value = 6.730088579440295
$ws.Range("C26").Value = 12.81526068145218
$ws.Range("C27").Value = 23.52963419638697
$ws.Range("C28").Value = 5.26052480300506
$ws.Range("C29").Value = 11.09573513214949
$ws.Range("C30").Value = 2.640079631252611
$ws.Range("C31").Value = 2.629507229983293
$ws.Range("C32").Value = 4.828566694001486
$ws.Range("C33").Value = 5.042280233945561
$ws.Range("C34").Value = 92.17774598110877
$ws.Range("C35").Value = 7.788839049410592
$ws.Range("C36").Value = 22.47994578464752
$ws.Range("C37").Value = 3.880071265839792
$ws.Range("C38").Value = 9.882174500735603
$ws.Range("C39").Value = 9.595209323425538
$ws.Range("C40").Value = 7.501873872100526
$ws.Range("C41").Value = 5.859375817760014
$ws.Range("C42").Value = 271.591759404743
